$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("subcorpus_bundle")

# Set the new / changed values in column E (order matters for shared-string table order)
$ws.Range("E22").Value = "inherited from count_bundle"
$ws.Range("E17").Value = "implemented explicitly for subcorpus_bundle"
$ws.Range("E12").Value = "implemented explicitly"
$ws.Range("E6").Value = "[not applicable]"
$ws.Range("E19").Value = "inherited from partition_bundle"
$ws.Range("E20").Value = "defined for subcorpus_bundle"
$ws.Range("E21").Value = "geerbt von partition_bundle"

# Update the active selection on the sheet
$ws.Range("E7").Select()
